$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.865.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.63%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.147.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.49%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "201.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.25%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "625.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.71%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.218"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.36%  "

# Row 9
$ws.Range("E9").Value = "  +1.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.486"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.39%  "

# Row 11
$ws.Range("E11").Value = "  +0.58%  "

# Row 12
$ws.Range("E12").Value = "  +6.99%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.720.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.25%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.57%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000204"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.33%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.723.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.49%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.144.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.56%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.97%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.92%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +20.49%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "404.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.34%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.78%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.21%  "

# Row 24
$ws.Range("B24").Value = "NEARProtocol"
$ws.Range("C24").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.41%  "

# Row 25
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.300.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.25%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "73.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.80%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.74%  "

# Row 28
$ws.Range("E28").Value = "  +0.35%  "

# Row 29
$ws.Range("E29").Value = "  +4.86%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.993"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.66%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.37%  "

# Row 32
$ws.Range("E32").Value = "  +5.38%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "523.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.28%  "

# Row 34
$ws.Range("E34").Value = "  +7.61%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.136"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +21.82%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "21.70"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.67%  "

# Row 37
$ws.Range("E37").Value = "  +0.03%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.42%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "196.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.64%  "

# Row 40
$ws.Range("E40").Value = "  -1.31%  "

# Row 41
$ws.Range("E41").Value = "  +0.59%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.103"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.79%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.54%  "

# Row 44
$ws.Range("E44").Value = "  +0.03%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.811"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +21.16%  "

# Row 46
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.90%  "

# Row 47
$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.26%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "42.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.19%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.45%  "

# Row 50
$ws.Range("E50").Value = "  +5.08%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.59%  "
